$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# The existing row 4 (a test-script row exercising "excluirEmpresa")
# needs to move down to row 5 to make room for a new "cadastrarEmpresa"
# test row (company "Kalunga"). Inserting a row above row 4 shifts the
# old row 4 down to row 5.
# ---------------------------------------------------------------------
$ws.Rows("4:4").Insert(-4121)

# The shifted-down row keeps the same test data, just with the next id.
$ws.Range("A5").Value = 4

# ---------------------------------------------------------------------
# Fill in the new row 4 with the "cadastrarEmpresa" test data (Kalunga).
# Values are set first, then each cell's formatting is re-applied by
# copying from a cell that already carries the exact desired style so
# the workbook ends up with the same visual formatting as before.
# ---------------------------------------------------------------------
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "formEmpresa"
$ws.Range("C4").Value = "cadastrarEmpresa"
$ws.Range("D4").Value = "formEmpresa"
$ws.Range("E4").Value = "Kalunga"

# Numeric-looking CNPJ value -- must stay text so the leading zeros
# survive (same trick the workbook already used for the other CNPJ in
# column F).
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = "00656601000157"

$ws.Range("G4").Value = "Kalunga loja"
$ws.Range("H4").Value = "Rua Augusta"
$ws.Range("I4").Value = 777788
$ws.Range("J4").Value = "silva"
$ws.Range("K4").Value = 666666
$ws.Range("L4").Value = "faturamento"
$ws.Range("M4").Value = "silva@gmail.com"
$ws.Range("N4").Value = "cadastro realizado com sucesso"

# Re-apply formatting (number format + quote-prefix flag) cell by cell
# from rows that already carry the right style.
$ws.Range("C3").Copy()
$ws.Range("C4").PasteSpecial(-4122)

$ws.Range("D3").Copy()
$ws.Range("D4").PasteSpecial(-4122)

$ws.Range("E3").Copy()
$ws.Range("E4").PasteSpecial(-4122)

$ws.Range("F3").Copy()
$ws.Range("F4").PasteSpecial(-4122)

$ws.Range("G3").Copy()
$ws.Range("G4").PasteSpecial(-4122)

$ws.Range("H3").Copy()
$ws.Range("H4").PasteSpecial(-4122)

$ws.Range("I3").Copy()
$ws.Range("I4").PasteSpecial(-4122)

# Column J ("coluna9", a person name) should use the same style as the
# other text name-style cells (e.g. D3), not the phone-style format
# that sits above it in J3.
$ws.Range("D3").Copy()
$ws.Range("J4").PasteSpecial(-4122)

$ws.Range("K3").Copy()
$ws.Range("K4").PasteSpecial(-4122)

$ws.Range("L3").Copy()
$ws.Range("L4").PasteSpecial(-4122)

$ws.Range("N3").Copy()
$ws.Range("N4").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Hyperlink the new email address cell, then restore the same cell
# formatting used by the other hyperlink cell (M3), since Hyperlinks.Add
# applies its own default hyperlink style.
# ---------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("M4"), "mailto:silva@gmail.com", [Type]::Missing, [Type]::Missing, "silva@gmail.com")
$ws.Range("M3").Copy()
$ws.Range("M4").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Selection reflects where editing left off.
# ---------------------------------------------------------------------
$ws.Range("N5").Select()
